# AFDP-835 - ARK-I - Close complaint/rework/resubmit tasks not showing on the
# task node in complaints module - grant read access by default to tasks
#
# 1) "Complaint - Restrict Access to Drafts" becomes "Complaint - Grant Access
#    to Drafts" and its action flips from "deny read to *" to "grant read to *"
# 2) A new rule row is appended for TASK objects that grants read access by
#    default.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the existing complaint drafts rule (row 21) - rename it and flip the
# access action from deny to grant.
$ws.Range("B21").Value = "Complaint – Grant Access to Drafts"
$ws.Range("G21").Value = "grant read to *"
$ws.Rows.Item(21).RowHeight = 23.5

# Add the new "Task - default read access" rule on the first empty row (30).
$ws.Range("B30").Value = "Task – default read access"
$ws.Range("C30").Value = "TASK"
$ws.Range("G30").Value = "grant read to *"
$ws.Rows.Item(30).RowHeight = 13.8

# Leave the cursor on the next empty row, matching where Excel would land
# after typing the new row's data.
$ws.Range("B31").Select() | Out-Null
